$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the two affected blocks as Text first so the new values (which look
# numeric/percent, e.g. "314.00" or "0.40%") are stored as literal strings,
# matching the inlineStr cells already in the workbook, not auto-converted
# to numbers/percentages by Excel's input parser.
$block1 = $ws.Range("D2:E27")
$block2 = $ws.Range("D39:E51")
$block1.NumberFormat = "@"
$block2.NumberFormat = "@"

$ws.Range("D2").Value = "314.00"
$ws.Range("E2").Value = "0.40%"
$ws.Range("D3").Value = "37.23"
$ws.Range("E3").Value = "-1.10%"
$ws.Range("D4").Value = "5.127"
$ws.Range("E4").Value = "-0.53%"
$ws.Range("D5").Value = "0.07919"
$ws.Range("E5").Value = "0.36%"
$ws.Range("D6").Value = "8.430"
$ws.Range("E6").Value = "1.57%"
$ws.Range("D7").Value = "1.887"
$ws.Range("E7").Value = "-1.54%"
$ws.Range("D8").Value = "2.976"
$ws.Range("E8").Value = "4.30%"
$ws.Range("D9").Value = "0.9323"
$ws.Range("E9").Value = "1.32%"
$ws.Range("D10").Value = "0.1267"
$ws.Range("E10").Value = "3.91%"
$ws.Range("D11").Value = "0.1922"
$ws.Range("E11").Value = "-0.05%"
$ws.Range("D12").Value = "0.08965"
$ws.Range("E12").Value = "-2.13%"
$ws.Range("D13").Value = "0.03347"
$ws.Range("E13").Value = "0.82%"
$ws.Range("D14").Value = "0.09519"
$ws.Range("E14").Value = "-0.62%"
$ws.Range("D15").Value = "0.001391"
$ws.Range("E15").Value = "0.03%"
$ws.Range("D16").Value = "0.006125"
$ws.Range("E16").Value = "7.91%"
$ws.Range("D17").Value = "3.391"
$ws.Range("E17").Value = "-3.42%"
$ws.Range("D18").Value = "4.431"
$ws.Range("E18").Value = "0.44%"
$ws.Range("D19").Value = "0.3491"
$ws.Range("E19").Value = "1.41%"
$ws.Range("D20").Value = "6.484"
$ws.Range("E20").Value = "23.32%"
$ws.Range("D21").Value = "0.1300"
$ws.Range("E21").Value = "2.14%"
$ws.Range("E22").Value = "-11.42%"
$ws.Range("D23").Value = "0.04346"
$ws.Range("E23").Value = "-0.32%"
$ws.Range("E24").Value = "-4.35%"
$ws.Range("D25").Value = "0.004235"
$ws.Range("E25").Value = "-10.25%"
$ws.Range("D26").Value = "0.0001324"
$ws.Range("E26").Value = "8.33%"
$ws.Range("D27").Value = "0.0003953"
$ws.Range("D39").Value = "0.02309"
$ws.Range("E39").Value = "0.46%"
$ws.Range("E40").Value = "0.59%"
$ws.Range("D41").Value = "0.007470"
$ws.Range("E41").Value = "0.00%"
$ws.Range("D42").Value = "0.1387"
$ws.Range("E42").Value = "2.03%"
$ws.Range("D43").Value = "0.008443"
$ws.Range("E43").Value = "-7.11%"
$ws.Range("D44").Value = "0.002061"
$ws.Range("E44").Value = "5.69%"
$ws.Range("D45").Value = "0.007931"
$ws.Range("E45").Value = "-7.76%"
$ws.Range("D46").Value = "0.00006314"
$ws.Range("E46").Value = "-4.60%"
$ws.Range("D47").Value = "0.00000000747"
$ws.Range("E47").Value = "-0.35%"
$ws.Range("D48").Value = "0.002854"
$ws.Range("E48").Value = "-14.71%"
$ws.Range("D49").Value = "0.001678"
$ws.Range("E49").Value = "39.74%"
$ws.Range("D50").Value = "0.00002091"
$ws.Range("E50").Value = "-0.35%"
$ws.Range("D51").Value = "0.0001992"
$ws.Range("E51").Value = "-0.35%"

# Restore the original (default/general) formatting now that the text values
# are locked in, so we don't leave a stray NumberFormat behind on these cells.
$block1.ClearFormats()
$block2.ClearFormats()
